# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New values for column G ("K") rows 2-27, replacing the old Strike# counts
$kValues = @{
    2  = 4
    3  = 4
    4  = 7
    5  = 9
    6  = 5
    7  = 11
    8  = 3
    9  = 8
    10 = 6
    11 = 11
    12 = 5
    13 = 6
    14 = 8
    15 = 11
    16 = 9
    17 = 11
    18 = 4
    19 = 6
    20 = 7
    21 = 9
    22 = 10
    23 = 10
    24 = 7
    25 = 6
    26 = 7
    27 = 5
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
